$wb = $excel.ActiveWorkbook
$wsConstructor = $wb.Worksheets.Item("Nutrients_Constructor")
$wsScale = $wb.Worksheets.Item("Nutrients_Scale")

# --- Add the two new test-case rows to the Nutrients_Scale sheet ---
# Typed in the same order a tester filling in the "smaller portion" / "bigger
# portion" scale() test cases would naturally use, so new shared strings land
# in the same dictionary order as the authored workbook.
$wsScale.Range("E3").Value = "Nutrients(100,45,10,1,0)"
$wsScale.Range("F3").Value = "Nutrients(50, 22.5, 5, 0.5, 0)"
$wsScale.Range("C3").Value = "Smaller portion"
$wsScale.Range("C4").Value = "Bigger portion"
$wsScale.Range("B3").Value = "Valid"
$wsScale.Range("B4").Value = "Valid"
$wsScale.Range("F4").Value = "Nutrients(200, 90, 20, 2, 0)"

$wsScale.Range("A3").Value = 1
$wsScale.Range("D3").Value = 50
$wsScale.Range("D4").Value = 200
$wsScale.Range("E4").Value = "Nutrients(100,45,10,1,0)"

# --- Widen the columns so the new, longer values stay readable ---
$wsScale.Range("B1:C1").ColumnWidth = 13.166666666666666
$wsScale.Range("E1").ColumnWidth = 20.5
$wsScale.Range("F1").ColumnWidth = 23.5

# --- Selection / active-sheet bookkeeping left behind by the editing session ---
$wsConstructor.Range("H12").Select() | Out-Null
$wsScale.Activate()
$wsScale.Range("F5").Select() | Out-Null
